$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2579a007-7ce0-45df-9b42-f43695ee8fe6_fila_4.png"
$ws.Range("B5").Value = "Yadder Fernando Torres"
$ws.Range("A6").Value = "2fa5f1fa-52da-410b-bb24-4b81feae4923_fila_9.png"
$ws.Range("B6").Value = "Roman Alfonso Grios Boza"
$ws.Range("A7").Value = "387f12a5-1bfc-4f1a-8266-e9bb5d3e65e8_fila_14.png"
$ws.Range("B7").Value = "Angel Isaac Alvarez Quiñonez"
$ws.Range("A8").Value = "3900499e-472d-46ac-bfc6-94423a039dea_fila_13.png"
$ws.Range("B8").Value = "Marlon Josue Gonzales Cano"
$ws.Range("A9").Value = "5fdfb4ee-51c6-4975-890d-de3fd11b4ae0_fila_5.png"
$ws.Range("B9").Value = "Erick Espinoza"
$ws.Range("A10").Value = "804e69c0-37bb-4675-a506-6cec1fcb206e_fila_1.png"
$ws.Range("B10").Value = "Hotep Antonio Ruiz Lezama"
$ws.Range("A11").Value = "8b1dca71-2ffd-49e8-ae06-f050609aee13_fila_2.png"
$ws.Range("B11").Value = "Isabella Dompe Estrada"
$ws.Range("A12").Value = "96a2ae22-6463-4021-a822-307dc50678bc_fila_10.png"
$ws.Range("B12").Value = "Abraham Silva Ampre"
$ws.Range("A13").Value = "9a9c70a0-158f-4e69-8dcc-a6518bb51ba2_fila_6.png"
$ws.Range("B13").Value = "Ronier Jose Rivera"
$ws.Range("A14").Value = "9c8d2642-f51f-409f-99dd-b3d8034b8fc2_fila_8.png"
$ws.Range("B14").Value = "David Orlando Mena Valverd"
$ws.Range("A15").Value = "a0b46958-419f-4cc6-9e71-d471a237eaf5_fila_15.png"
$ws.Range("B15").Value = "Orlando Mauricio Guevara"
$ws.Range("A16").Value = "c630cac0-49a2-4c52-837c-96b9495fa9ab_fila_3.png"
$ws.Range("B16").Value = "Bryan Alexander Cano"
$ws.Range("A17").Value = "c9dcc57c-96d6-4867-9724-cf96190a2f89_fila_7.png"
$ws.Range("B17").Value = "Cristina Jozabed Carvajal"
$ws.Range("A18").Value = "e9016cc4-85a0-442f-9865-bda3a3086566_fila_11.png"
$ws.Range("B18").Value = "Eduardo Domingo Zeledon Merca"
$ws.Range("A19").Value = "e9071e59-145b-4d3b-853e-043cd79ae634_fila_12.png"
$ws.Range("B19").Value = "José Danilo Suárez"
